$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$ws1.Range("F2").Value = "2021-10-05 14:20:07.963575"
$ws1.Range("F3").Value = "2021-10-05 14:20:07.963583"
$ws1.Range("F4").Value = "2021-10-05 14:20:07.963586"
$ws1.Range("F5").Value = "2021-10-05 14:20:07.963589"
$ws1.Range("F6").Value = "2021-10-05 14:20:07.963592"
$ws1.Range("F7").Value = "2021-10-05 14:20:07.963594"
$ws1.Range("F8").Value = "2021-10-05 14:20:07.963597"
$ws1.Range("F9").Value = "2021-10-05 14:20:07.963599"
$ws1.Range("F10").Value = "2021-10-05 14:20:07.963602"
$ws1.Range("F11").Value = "2021-10-05 14:20:07.963605"
$ws1.Range("F12").Value = "2021-10-05 14:20:07.963608"
$ws1.Range("F13").Value = "2021-10-05 14:20:07.963610"
$ws1.Range("F14").Value = "2021-10-05 14:20:07.963613"
$ws1.Range("F15").Value = "2021-10-05 14:20:07.963615"
$ws1.Range("F16").Value = "2021-10-05 14:20:07.963618"
$ws1.Range("F17").Value = "2021-10-05 14:20:07.963620"
$ws1.Range("F18").Value = "2021-10-05 14:20:07.963623"
$ws1.Range("F19").Value = "2021-10-05 14:20:07.963626"
$ws1.Range("F20").Value = "2021-10-05 14:20:07.963628"
$ws1.Range("F21").Value = "2021-10-05 14:20:07.963631"
$ws1.Range("F22").Value = "2021-10-05 14:20:07.963633"

# --- Add the new "metadata" sheet, placed after "data" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "metadata"

# --- Header row (row 1), columns B:G ---
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# --- Data row (row 2) ---
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Epidermolysis bullosa"
$newSheet.Range("C2").Value = 119
# data_version "1.6" must be stored as TEXT, not a number: force a text
# number format while assigning, then strip the format again so the cell
# ends up with the workbook's default (unstyled) appearance but a string value.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.6"
$newSheet.Range("D2").ClearFormats()
$newSheet.Range("E2").Value = "2019-01-07T16:40:38.182703Z"
$newSheet.Range("F2").Value = "2021-10-05 14:20:07.959857"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/119/?format=json"

# --- Apply the same header style used on the "data" sheet (bold, centered, bordered) ---
$ws1.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$null = $newSheet.Range("A1").Select()
